# Fruta / hortaliza, semanal
# Insert two new rows of data (rows 27 and 28) above the existing row 27,
# pushing the former rows 27-55 down to 29-57, and populate the two new
# rows with the new "Angeleno" price records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 27. Doing this twice at the
# same index shifts the prior row 27 (and everything below it) down by two.
$ws.Rows.Item(27).Insert()
$ws.Rows.Item(27).Insert()

# ---- New row 27 ----
$ws.Range("A27").Value = 7
$ws.Range("B27").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C27").Value = "Ñuble"
$ws.Range("D27").Value = 44645
$ws.Range("E27").Value = 16
$ws.Range("F27").Value = "Fruta"
$ws.Range("G27").Value = 100103
$ws.Range("H27").Value = "Frutos de hueso (carozo)"
$ws.Range("I27").Value = 100103002
$ws.Range("J27").Value = "Ciruela"
$ws.Range("K27").Value = "Angeleno"
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 120
$ws.Range("N27").Value = 8000
$ws.Range("O27").Value = 9000
$ws.Range("P27").Value = 8500
$ws.Range("Q27").Value = "$/bandeja 18 kilos granel"
$ws.Range("R27").Value = "Provincia de Curicó"
$ws.Range("S27").Value = 472
$ws.Range("T27").Value = 18

# ---- New row 28 ----
$ws.Range("A28").Value = 7
$ws.Range("B28").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C28").Value = "Ñuble"
$ws.Range("D28").Value = 44645
$ws.Range("E28").Value = 16
$ws.Range("F28").Value = "Fruta"
$ws.Range("G28").Value = 100103
$ws.Range("H28").Value = "Frutos de hueso (carozo)"
$ws.Range("I28").Value = 100103002
$ws.Range("J28").Value = "Ciruela"
$ws.Range("K28").Value = "Angeleno"
$ws.Range("L28").Value = "Segunda"
$ws.Range("M28").Value = 60
$ws.Range("N28").Value = 7000
$ws.Range("O28").Value = 7000
$ws.Range("P28").Value = 7000
$ws.Range("Q28").Value = "$/bandeja 18 kilos granel"
$ws.Range("R28").Value = "Provincia de Curicó"
$ws.Range("S28").Value = 389
$ws.Range("T28").Value = 18
